$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text for the "R10" rule row
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active cell selection left by the editor
$ws.Range("E8").Select()
